$p = $ppt.ActivePresentation

# --- 1) Refresh the cached "datetimeFigureOut" date placeholder text ---
#     (slide master, every slide layout, and the notes master) from
#     1/30/2019 to 6/24/2019, as happens when PowerPoint re-saves the
#     deck on a later day and re-caches the auto date fields.
$newDate = "6/24/2019"

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes $newDate

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Set-DatePlaceholderText $layouts.Item($L).Shapes $newDate
}

# Notes master (its date placeholder only updates via the HeadersFooters API)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# --- 2) Update the subtitle on slide 1: "Winter 2019" -> "2019" ---
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "Winter 2019") {
        $shp.TextFrame.TextRange.Text = "2019"
    }
}
